$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume snapshot refresh.
# Cells B2:E51 are stored as literal text in the source sheet (prices use
# dotted thousands separators like '69.505.53' and volumes carry padding
# spaces, e.g. '  +0.25%  '). Where a new value would otherwise be
# auto-recognised as a number (e.g. '654.17', '1.00'), force the cell to
# Text format first so it is written back as text, not coerced to a number.

$ws.Range("D2").Value = '69.508.88'

$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '3.676.67'

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '654.17'

$ws.Range("E5").Value = '  -3.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.59'

$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("E8").Value = '  +0.36%  '

$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("E10").Value = '  -0.37%  '

$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").Value = '4.295.67'

$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.52'

$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("D15").Value = '3.683.03'

$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("D16").Value = '69.502.52'

$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("E18").Value = '  -0.32%  '

$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.47'

$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.69'

$ws.Range("E21").Value = '  -2.83%  '

$ws.Range("E22").Value = '  -1.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.58'

$ws.Range("E23").Value = '  -0.38%  '

$ws.Range("D24").Value = '3.823.77'

$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("E26").Value = '  +1.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.85'

$ws.Range("E27").Value = '  -0.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.92'

$ws.Range("E28").Value = '  -2.12%  '

$ws.Range("E30").Value = '  -5.51%  '

$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("E32").Value = '  +0.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.62'

$ws.Range("E33").Value = '  -1.24%  '

$ws.Range("E34").Value = '  -2.68%  '

$ws.Range("B35").Value = 'RenzoRestakedETH'

$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'

$ws.Range("D35").Value = '3.667.14'

$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("B36").Value = 'Kaspa'

$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.162'

$ws.Range("E36").Value = '  +2.70%  '

$ws.Range("E37").Value = '  +1.38%  '

$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.98'

$ws.Range("E39").Value = '  -3.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '178.99'

$ws.Range("E40").Value = '  +5.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("E42").Value = '  -1.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0895'

$ws.Range("E43").Value = '  -1.18%  '

$ws.Range("E44").Value = '  -1.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.86'

$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.27'

$ws.Range("E47").Value = '  -0.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.19'

$ws.Range("E48").Value = '  -3.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000266'

$ws.Range("E49").Value = '  -4.85%  '

$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("E51").Value = '  -3.81%  '
